$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$col = $ws.Columns.Item(2)

$col.Replace("IcrisatBW5CvM35-1FertHighIrrigHigh", "BW5CvM35-1FertHighIrrigHigh", 1, 1, $false, $false, $true) | Out-Null
$col.Replace("IcrisatBW5CvCSH13RFertHighIrrigHigh", "BW5CvCSH13RFertHighIrrigHigh", 1, 1, $false, $false, $true) | Out-Null
$col.Replace("IcrisatBW5CvAtx623xRTX430FertHighIrrigHigh", "BW5CvAtx623xRTX430FertHighIrrigHigh", 1, 1, $false, $false, $true) | Out-Null
$col.Replace("IcrisatBW5CvQL41xQL36FertHighIrrigHigh", "BW5CvQL41xQL36FertHighIrrigHigh", 1, 1, $false, $false, $true) | Out-Null
$col.Replace("IcrisatBW5CvM35-1FertMedIrrigLow", "BW5CvM35-1FertMedIrrigLow", 1, 1, $false, $false, $true) | Out-Null
$col.Replace("IcrisatBW5CvCSH13RFertMedIrrigLow", "BW5CvCSH13RFertMedIrrigLow", 1, 1, $false, $false, $true) | Out-Null
$col.Replace("IcrisatBW5CvAtx623xRTX430FertMedIrrigLow", "BW5CvAtx623xRTX430FertMedIrrigLow", 1, 1, $false, $false, $true) | Out-Null
$col.Replace("IcrisatBW5CvQL41xQL36FertMedIrrigLow", "BW5CvQL41xQL36FertMedIrrigLow", 1, 1, $false, $false, $true) | Out-Null
$col.Replace("IcrisatBW5CvM35-1FertLowIrrigLow", "BW5CvM35-1FertLowIrrigLow", 1, 1, $false, $false, $true) | Out-Null
$col.Replace("IcrisatBW5CvCSH13RFertLowIrrigLow", "BW5CvCSH13RFertLowIrrigLow", 1, 1, $false, $false, $true) | Out-Null
$col.Replace("IcrisatBW5CvAtx623xRTX430FertLowIrrigLow", "BW5CvAtx623xRTX430FertLowIrrigLow", 1, 1, $false, $false, $true) | Out-Null
$col.Replace("IcrisatBW5CvQL41xQL36FertLowIrrigLow", "BW5CvQL41xQL36FertLowIrrigLow", 1, 1, $false, $false, $true) | Out-Null
$col.Replace("IcrisatBW8CvCSH13RFertHighIrrigOn", "BW8CvCSH13RFertHighIrrigOn", 1, 1, $false, $false, $true) | Out-Null
$col.Replace("IcrisatBW8CvM35-1FertHighIrrigOn", "BW8CvM35-1FertHighIrrigOn", 1, 1, $false, $false, $true) | Out-Null
$col.Replace("IcrisatBW8CvAtx623xRTx430FertHighIrrigOn", "BW8CvAtx623xRTx430FertHighIrrigOn", 1, 1, $false, $false, $true) | Out-Null
$col.Replace("IcrisatBW8CvQL41xQL36FertHighIrrigOn", "BW8CvQL41xQL36FertHighIrrigOn", 1, 1, $false, $false, $true) | Out-Null
$col.Replace("IcrisatBW8CvM35-1FertMedIrrigOff", "BW8CvM35-1FertMedIrrigOff", 1, 1, $false, $false, $true) | Out-Null
$col.Replace("IcrisatBW8CvCSH13RFertMedIrrigOff", "BW8CvCSH13RFertMedIrrigOff", 1, 1, $false, $false, $true) | Out-Null
$col.Replace("IcrisatBW8CvATX623xRTX430FertMedIrrigOff", "BW8CvATX623xRTX430FertMedIrrigOff", 1, 1, $false, $false, $true) | Out-Null
$col.Replace("IcrisatBW8CvQL41xQL36FertMedIrrigOff", "BW8CvQL41xQL36FertMedIrrigOff", 1, 1, $false, $false, $true) | Out-Null
$col.Replace("IcrisatBW8CvM35-1FertLowIrrigOff", "BW8CvM35-1FertLowIrrigOff", 1, 1, $false, $false, $true) | Out-Null
$col.Replace("IcrisatBW8CvCSH13RFertLowIrrigOff", "BW8CvCSH13RFertLowIrrigOff", 1, 1, $false, $false, $true) | Out-Null
$col.Replace("IcrisatBW8CvATX623xRTX430FertLowIrrigOff", "BW8CvATX623xRTX430FertLowIrrigOff", 1, 1, $false, $false, $true) | Out-Null
$col.Replace("IcrisatBW8CvQL41xQL36FertLowIrrigOff", "BW8CvQL41xQL36FertLowIrrigOff", 1, 1, $false, $false, $true) | Out-Null
$col.Replace("Hermitage1996FertHighIrrigOnCvBuster", "HE1-4FertHighIrrigOnCvBuster", 1, 1, $false, $false, $true) | Out-Null
$col.Replace("Hermitage1996FertHighIrrigOnCvQL41xQL36", "HE1-4FertHighIrrigOnCvQL41xQL36", 1, 1, $false, $false, $true) | Out-Null
$col.Replace("Hermitage1996FertHighIrrigOnCvM351", "HE1-4FertHighIrrigOnCvM351", 1, 1, $false, $false, $true) | Out-Null
$col.Replace("Hermitage1996FertLowIrrigOnCvQL41xQL36", "HE1-4FertLowIrrigOnCvQL41xQL36", 1, 1, $false, $false, $true) | Out-Null
$col.Replace("Hermitage1996FertLowIrrigOnCvBuster", "HE1-4FertLowIrrigOnCvBuster", 1, 1, $false, $false, $true) | Out-Null
$col.Replace("Hermitage1996FertLowIrrigOnCvM351", "HE1-4FertLowIrrigOnCvM351", 1, 1, $false, $false, $true) | Out-Null
$col.Replace("Hermitage1996FertHighIrrigOffCvQL41xQL36", "HE1-4FertHighIrrigOffCvQL41xQL36", 1, 1, $false, $false, $true) | Out-Null
$col.Replace("Hermitage1996FertHighIrrigOffCvBuster", "HE1-4FertHighIrrigOffCvBuster", 1, 1, $false, $false, $true) | Out-Null
$col.Replace("Hermitage1996FertHighIrrigOffCvM351", "HE1-4FertHighIrrigOffCvM351", 1, 1, $false, $false, $true) | Out-Null
$col.Replace("Hermitage1996FertLowIrrigOffCvQL41xQL36", "HE1-4FertLowIrrigOffCvQL41xQL36", 1, 1, $false, $false, $true) | Out-Null
$col.Replace("Hermitage1996FertLowIrrigOffCvBuster", "HE1-4FertLowIrrigOffCvBuster", 1, 1, $false, $false, $true) | Out-Null
$col.Replace("Hermitage1996FertLowIrrigOffCvM351", "HE1-4FertLowIrrigOffCvM351", 1, 1, $false, $false, $true) | Out-Null
$col.Replace("Hermitage1997FertHighIrrigHighCvBuster", "HE5-8FertHighIrrigHighCvBuster", 1, 1, $false, $false, $true) | Out-Null
$col.Replace("Hermitage1997FertHighIrrigHighCvM351", "HE5-8FertHighIrrigHighCvM351", 1, 1, $false, $false, $true) | Out-Null
$col.Replace("Hermitage1997FertHighIrrigHighCvCSH13R", "HE5-8FertHighIrrigHighCvCSH13R", 1, 1, $false, $false, $true) | Out-Null
$col.Replace("Hermitage1997FertLowIrrigHighCvBuster", "HE5-8FertLowIrrigHighCvBuster", 1, 1, $false, $false, $true) | Out-Null
$col.Replace("Hermitage1997FertLowIrrigHighCvM351", "HE5-8FertLowIrrigHighCvM351", 1, 1, $false, $false, $true) | Out-Null
$col.Replace("Hermitage1997FertLowIrrigHighCvCSH13R", "HE5-8FertLowIrrigHighCvCSH13R", 1, 1, $false, $false, $true) | Out-Null
$col.Replace("Hermitage1997FertHighIrrigLowCvBuster", "HE5-8FertHighIrrigLowCvBuster", 1, 1, $false, $false, $true) | Out-Null
$col.Replace("Hermitage1997FertHighIrrigLowCvM351", "HE5-8FertHighIrrigLowCvM351", 1, 1, $false, $false, $true) | Out-Null
$col.Replace("Hermitage1997FertHighIrrigLowCvCSH13R", "HE5-8FertHighIrrigLowCvCSH13R", 1, 1, $false, $false, $true) | Out-Null
$col.Replace("Hermitage1997FertLowIrrigLowCvBuster", "HE5-8FertLowIrrigLowCvBuster", 1, 1, $false, $false, $true) | Out-Null
$col.Replace("Hermitage1997FertLowIrrigLowCvM351", "HE5-8FertLowIrrigLowCvM351", 1, 1, $false, $false, $true) | Out-Null
$col.Replace("Hermitage1997FertLowIrrigLowCvCSH13R", "HE5-8FertLowIrrigLowCvCSH13R", 1, 1, $false, $false, $true) | Out-Null
$col.Replace("Lawes1995FertLowirrigOffCvBuster", "LE13FertLowirrigOffCvBuster", 1, 1, $false, $false, $true) | Out-Null
$col.Replace("Lawes1995FertHighIrrigOffCvBuster", "LE13FertHighIrrigOffCvBuster", 1, 1, $false, $false, $true) | Out-Null
$col.Replace("Lawes1995FertHighIrrigOnCvBuster", "LE13FertHighIrrigOnCvBuster", 1, 1, $false, $false, $true) | Out-Null
$col.Replace("Lawes1995FertLowIrrigOffCvM351", "LE13FertLowIrrigOffCvM351", 1, 1, $false, $false, $true) | Out-Null
$col.Replace("Lawes1995FertHighIrrigOffCvM351", "LE13FertHighIrrigOffCvM351", 1, 1, $false, $false, $true) | Out-Null
$col.Replace("Lawes1995FertHighIrrigOnCvM351", "LE13FertHighIrrigOnCvM351", 1, 1, $false, $false, $true) | Out-Null
$col.Replace("Lawes1996EarlyCvBuster", "LE14CvBuster", 1, 1, $false, $false, $true) | Out-Null
$col.Replace("Lawes1996EarlyCvQL41xQL36", "LE14CvQL41xQL36", 1, 1, $false, $false, $true) | Out-Null
$col.Replace("Lawes1996EarlyCvM351", "LE14CvM351", 1, 1, $false, $false, $true) | Out-Null
$col.Replace("Lawes1996LateCvBuster", "LE15CvBuster", 1, 1, $false, $false, $true) | Out-Null
$col.Replace("Lawes1996LateCvQL41xQL36", "LE15CvQL41xQL36", 1, 1, $false, $false, $true) | Out-Null
$col.Replace("Lawes1996LateCvM351", "LE15CvM351", 1, 1, $false, $false, $true) | Out-Null
$col.Replace("Lawes1997LateCvBuster", "LE17CvBuster", 1, 1, $false, $false, $true) | Out-Null
$col.Replace("Lawes1997LateCvCSH13R", "LE17CvCSH13R", 1, 1, $false, $false, $true) | Out-Null
$col.Replace("Lawes1997LateCvM351", "LE17CvM351", 1, 1, $false, $false, $true) | Out-Null
$col.Replace("Lawes1998FertMedCvBuster", "LE19FertMedCvBuster", 1, 1, $false, $false, $true) | Out-Null
$col.Replace("Lawes1998FertLowCvCSH13R", "LE19FertLowCvCSH13R", 1, 1, $false, $false, $true) | Out-Null
$col.Replace("Lawes1998FertOffCvBuster", "LE19FertOffCvBuster", 1, 1, $false, $false, $true) | Out-Null
$col.Replace("Lawes1998FertLowCvBuster", "LE19FertLowCvBuster", 1, 1, $false, $false, $true) | Out-Null
$col.Replace("Lawes1998FertHighCvBuster", "LE19FertHighCvBuster", 1, 1, $false, $false, $true) | Out-Null
$col.Replace("Lawes1998FertOffCvCSH13R", "LE19FertOffCvCSH13R", 1, 1, $false, $false, $true) | Out-Null
$col.Replace("Lawes1998FertMedCvCSH13R", "LE19FertMedCvCSH13R", 1, 1, $false, $false, $true) | Out-Null
$col.Replace("Lawes1998FertHighCvCSH13R", "LE19FertHighCvCSH13R", 1, 1, $false, $false, $true) | Out-Null
$col.Replace("Lawes1999CvA35xQL36FertLow", "LE21CvA35xQL36FertLow", 1, 1, $false, $false, $true) | Out-Null
$col.Replace("Lawes1999CvA35xQL36FertHigh", "LE21CvA35xQL36FertHigh", 1, 1, $false, $false, $true) | Out-Null
$col.Replace("Lawes1999CvCSH13RFertLow", "LE21CvCSH13RFertLow", 1, 1, $false, $false, $true) | Out-Null
$col.Replace("Lawes1999CvQL39xQL36FertLow", "LE21CvQL39xQL36FertLow", 1, 1, $false, $false, $true) | Out-Null
$col.Replace("Lawes1999CvCSH13RFertMed", "LE21CvCSH13RFertMed", 1, 1, $false, $false, $true) | Out-Null
$col.Replace("Lawes1999CvA35xQL36FertMed", "LE21CvA35xQL36FertMed", 1, 1, $false, $false, $true) | Out-Null
$col.Replace("Lawes1999CvQL39xQL36FertMed", "LE21CvQL39xQL36FertMed", 1, 1, $false, $false, $true) | Out-Null
$col.Replace("Lawes1999CvQL39xQL36FertHigh", "LE21CvQL39xQL36FertHigh", 1, 1, $false, $false, $true) | Out-Null
$col.Replace("Lawes1999CvCSH13RFertHigh", "LE21CvCSH13RFertHigh", 1, 1, $false, $false, $true) | Out-Null

# Update row heights for rows 14 and 15 (15 -> 13.8 points)
$ws.Rows.Item(14).RowHeight = 13.8
$ws.Rows.Item(15).RowHeight = 13.8

# Update view/selection state to match the saved workbook view
$ws.Range("B982").Select()
